# B6-PowerPoint.pptx edit script
# 1) Re-style the three data tables (slides 14, 15, 16) to the
#    "{0B55327C-AB88-4A1A-B1A3-79D772C16EF8}" built-in table style.
# 2) Swap the presentation's active theme colour scheme back to the
#    stock "Office" palette (it was previously the "Integral / Red
#    Violet" palette).

$p = $ppt.ActivePresentation

$newTableStyle = "{0B55327C-AB88-4A1A-B1A3-79D772C16EF8}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# Restore the default Office colour scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) on the deck's theme.
$themeColors = @{
    1  = 0        # dk1      000000
    2  = 16777215 # lt1      FFFFFF
    3  = 6968388  # dk2      44546A
    4  = 15132391 # lt2      E7E6E6
    5  = 13998939 # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

$scheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $themeColors[$i]
}
